# Fix Training Data Issue
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") holds the game date as text "2-15-2012-13"; correct
# it to ISO-style "2013-02-15" for every data row (rows 2-31).
#
# The date text must stay literal TEXT (not get auto-converted into a
# date serial by Excel's smart-typing), so the range is marked as Text
# ("@") before the corrected values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2013-02-15"
}
